# Insert one new weekly price record for "Espinaca" (Vega Modelo de Temuco)
# as row 229, pushing the existing rows 229:288 down to 230:289.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 229 down by one row (classic Excel "insert row").
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with the new record's data.
$ws.Cells.Item(229, 1).Value = 10
$ws.Cells.Item(229, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(229, 3).Value = "La Araucanía"
$ws.Cells.Item(229, 4).Value = 45120
$ws.Cells.Item(229, 5).Value = 9
$ws.Cells.Item(229, 6).Value = 100112012
$ws.Cells.Item(229, 7).Value = "Espinaca"
$ws.Cells.Item(229, 8).Value = "Sin especificar"
$ws.Cells.Item(229, 9).Value = "Primera"
$ws.Cells.Item(229, 10).Value = 80
$ws.Cells.Item(229, 11).Value = 8000
$ws.Cells.Item(229, 12).Value = 8000
$ws.Cells.Item(229, 13).Value = 8000
$ws.Cells.Item(229, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(229, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(229, 16).Value = 667
$ws.Cells.Item(229, 17).Value = 12
$ws.Cells.Item(229, 18).Value = "Hortaliza"
